$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6, 18, 32, 71 were touched (re-pasted) by the original edit: values are
# unchanged but the cells now carry an explicit "Normal" style and previously
# blank trailing cells (J/S/T) materialize as empty cells.
$ws.Range("A6:T6").Style = "Normal"
$ws.Range("A18:T18").Style = "Normal"
$ws.Range("A32:T32").Style = "Normal"
$ws.Range("A71:T71").Style = "Normal"

# New row 84: OdtPainting
$ws.Range("A84").Value = "OdtPainting"
$ws.Range("B84").Value = "An experiment at ODT stage."
$ws.Range("C84").Value = "TOP"
$ws.Range("D84").Value = "EvapDOdt1"
$ws.Range("E84").Value = "None"
$ws.Range("F84").Value = 4
$ws.Range("G84").Value = "Painting"
$ws.Range("H84").Value = "AtomNumber;DensityFit;CenterFit"
$ws.Range("I84").Value = "LSR"
$ws.Range("J84").Value = "[587 1430;820 1077;1185 1116;1412 1467;1386 1938;983 2154;629 1963;523 1643;509 1643]"
$ws.Range("K84").Value = "LF"
$ws.Range("L84").Value = "Absorption"
$ws.Range("M84").Value = "RandomPolarization"
$ws.Range("N84").Value = 8
$ws.Range("O84").Value = "BosonicGaussianFit1D"
$ws.Range("P84").Value = 15
$ws.Range("Q84").Value = "LinearFit1D"
$ws.Range("R84").Value = "StdErr"
$ws.Range("A84:T84").Style = "Normal"

# New row 85: NiBecTau
$ws.Range("A85").Value = "NiBecTau"
$ws.Range("B85").Value = "An experiment at the non-interacting BEC stage. Scan tau."
$ws.Range("C85").Value = "TOP"
$ws.Range("D85").Value = "Bec"
$ws.Range("E85").Value = "None"
$ws.Range("F85").Value = 4
$ws.Range("G85").Value = "tau"
$ws.Range("H85").Value = "DensityFit;AtomNumber;CenterFit"
$ws.Range("I85").Value = "LSR"
$ws.Range("J85").Value = "[883 1331;920 1373;977 1409;1044 1385;1079 1289;1076 1226;1060 1180;1016 1173;953 1169;903 1189;878 1246]"
$ws.Range("K85").Value = "NI"
$ws.Range("L85").Value = "Absorption"
$ws.Range("M85").Value = "StrongLight"
$ws.Range("N85").Value = 8
$ws.Range("O85").Value = "BosonicGaussianFit1D"
$ws.Range("P85").Value = 0.2
$ws.Range("Q85").Value = "ParabolicFit1D"
$ws.Range("R85").Value = "StdErr"
$ws.Range("T85").Value = "WaveformGeneratorName,XvWingMod;ChannelName,Ch1;WaveformListName,XvWingNi"
$ws.Range("A85:T85").Style = "Normal"

# New row 86: NiBecTofCameraOdt
$ws.Range("A86").Value = "NiBecTofCameraOdt"
$ws.Range("B86").Value = "A TOF experiment at the non-interacting BEC stage."
$ws.Range("C86").Value = "ODT"
$ws.Range("D86").Value = "SideOdtCamera"
$ws.Range("F86").Value = 4
$ws.Range("G86").Value = "TOF"
$ws.Range("H86").Value = "DensityFit;AtomNumber;CenterFit"
$ws.Range("I86").Value = "LSR"
$ws.Range("J86").Value = "[907 1249;924 1587;967 1593;1073 1593;1073 1468;1067 1326;1060 1215;1017 1158;918 1162]"
$ws.Range("K86").Value = "HF"
$ws.Range("L86").Value = "Absorption"
$ws.Range("M86").Value = "StrongLight"
$ws.Range("N86").Value = 8
$ws.Range("O86").Value = "BosonicGaussianFit1D"
$ws.Range("P86").Value = 1
$ws.Range("Q86").Value = "ParabolicFit1D"
$ws.Range("R86").Value = "StdErr"
$ws.Range("T86").Value = "WaveformGeneratorName,XvWingMod;ChannelName,Ch1;WaveformListName,XvWingNi"
$ws.Range("A86:T86").Style = "Normal"

$wb.Save()
